# Updated symbol list on Fri Jan 20 11:28:45 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for the crypto
# ticker rows on Sheet1. The source cells are stored as plain text
# (e.g. "287.89", "-0.81%"), so a leading apostrophe forces the COM
# layer to keep the new value as text instead of re-parsing it as a
# number/percentage; resetting Style back to "Normal" afterwards drops
# the transient quote-prefix formatting so no visible style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value = "'287.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.72%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.61%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'-0.63%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07315"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.39%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.338"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'27.89%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.732"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.67%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.722"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.21%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9026"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.65%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'17.91%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1693"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.25%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08277"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.14%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03123"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.82%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09943"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.66%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001493"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.72%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005827"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.13%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.494"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.71%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.097"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.01%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3329"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.26%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.13%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.192"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.71%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-12.17%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04500"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.07%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.33%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-9.84%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.15%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003394"
$ws.Range("D27").Style = "Normal"
$ws.Range("D39").Value = "'0.01579"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.10%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04447"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.91%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007342"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.40%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.87%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1326"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.02%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002221"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'10.15%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-5.53%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006116"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.58%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.561"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'10.87%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.13%"
$ws.Range("E51").Style = "Normal"
